$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Layout")

# Fill in the previously empty cells in rows 4 and 5 (F/G columns)
$ws.Range("F4").Value = "Minimum burst size - ttk.Label"
$ws.Range("G4").Value = "ttk.Entry"
$ws.Range("F5").Value = "Plot burst probability - ttk.Label"
$ws.Range("G5").Value = "checkbox"

# Update the active sheet's view/selection to G6 (and drop the old frozen top-left cell)
$ws.Activate()
$ws.Range("G6").Select()
